$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 7: date 2016-09-14 (42627), start 08:50 (0.368055...), end 22:00 (0.916666...)
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").NumberFormat = "h:mm"
$ws.Range("C7").NumberFormat = "h:mm"
$ws.Range("A7").Value = 42627
$ws.Range("B7").Value = 0.36805555555555558
$ws.Range("C7").Value = 0.91666666666666663
$ws.Range("D7").Value = "Implementação do login com ramificação. Implementação do CRUD de motoboy"

# Row 8: date 2016-09-14 (42627), start 22:00, end 22:45 (0.947916...)
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").NumberFormat = "h:mm"
$ws.Range("C8").NumberFormat = "h:mm"
$ws.Range("A8").Value = 42627
$ws.Range("B8").Value = 0.91666666666666663
$ws.Range("C8").Value = 0.94791666666666663
$ws.Range("D8").Value = "Estudo e tentativa de confirmação para exclusão da conta do cliente"

# Update selection to D9 (matches the diff's sheetView selection change)
$ws.Range("D9").Select()
